$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functions")

# Give the new "Symbol" column (C) a bit more room, matching the added
# C-column content below (round(a, accuracy) / modPow(base, exp, mod) / etc.)
$ws.Columns.Item(3).ColumnWidth = 20.877604166666668

# --- Status clean-up pass ---------------------------------------------

# POWER (row 8): was "IP", now "N/A"
$ws.Range("D8").Value = "N/A"

# MIN (row 15): add its missing symbol and mark it "Done"
$ws.Range("C15").Value = "min(a1, a2, …, aN)"
$ws.Range("D15").Value = "Done"

# FLOOR (row 18): was "N/A", now "IP"
$ws.Range("D18").Value = "IP"

# CEIL (row 19): was "N/A", now "IP"
$ws.Range("D19").Value = "IP"

# LN (row 20): was "IP", now "N/A"
$ws.Range("D20").Value = "N/A"

# LG (row 21): was "IP", now "N/A"
$ws.Range("D21").Value = "N/A"

# LOG (row 22): was "IP", now "N/A"
$ws.Range("D22").Value = "N/A"

# LOGN (row 23): was "IP", now "N/A"
$ws.Range("D23").Value = "N/A"

# RNG (row 24): was "N/A", now "Done"
$ws.Range("D24").Value = "Done"

# RNG_FLOAT (row 25): was "N/A", now "IP"
$ws.Range("D25").Value = "IP"

# --- New functions: ROUND and MOD_POWER ---------------------------------

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "ROUND"
$ws.Range("C44").Value = "round(a, accuracy)"
$ws.Range("D44").Value = "N/A"

$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "MOD_POWER"
$ws.Range("C45").Value = "modPow(base, exp, mod)"
$ws.Range("D45").Value = "IP"

# Leave the view pointed at the newly added row, like the author did.
$ws.Range("A45:XFD45").Select() | Out-Null
